$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Munka1")

# --- Row 16: test description text rephrased + hours updated ---
# (set first so the new shared string "elemi és bővített fgv.ek
# tesztelése" is registered before "Esztétika", matching the order
# new strings were appended in the original edit)
$ws.Range("B16").Value = "elemi és bővített fgv.ek tesztelése"
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 9

# --- Row 13: "Funkció" changes from "Adatbáziskezelés" to "Esztétika" ---
$ws.Range("A13").Value = "Esztétika"

# --- Move the closing "forrás" hyperlink row out of the way before
#     inserting a new row, since hyperlink anchors are not retargeted
#     automatically when rows shift. ---
$ws.Range("A19").Hyperlinks.Delete()
$ws.Range("A19").Clear()

# --- Insert a new row 18 (new task "Osztályok szeparálása" under the
#     "Esztétika" function), copying the formatting of row 17. ---
$ws.Rows.Item(18).Insert()
$ws.Range("A17:H17").Copy()
$ws.Range("A18:H18").PasteSpecial(-4104)
$ws.Range("A18").Value = "Esztétika"
$ws.Range("B18").Value = "Osztályok szeparálása"
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 5
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = "Pictori"

# --- Re-create the closing hyperlink row further down at row 27 ---
$ws.Range("A27").Value = "http://hungarian.joelonsoftware.com/Articles/PainlessSoftwareSchedules.html"
$ws.Range("A27").Style = "Hivatkozás"
$ws.Hyperlinks.Add($ws.Range("A27"), "http://hungarian.joelonsoftware.com/Articles/PainlessSoftwareSchedules.html")

# --- Update the worksheet's selected / active cell ---
$ws.Range("E22").Select()
